# Apply updated NATMI recomputation values (Efna1-Epha1) per Dr Hou's advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 23.60223533333333
$ws.Range("H2").Value = 70.80670599999999
$ws.Range("I2").Value = 0.8824726436021215
$ws.Range("J2").Value = 0.8824726436021214
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.497699333333333
$ws.Range("N2").Value = 13.493098
$ws.Range("O2").Value = 0.2265545985397048
$ws.Range("P2").Value = 0.2265545985397048
$ws.Range("Q2").Value = 106.1557581239098
$ws.Range("R2").Value = 955.4018231151879
$ws.Range("S2").Value = 0.1999282354935506
$ws.Range("T2").Value = 0.1999282354935506
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 23.60223533333333
$ws.Range("H3").Value = 70.80670599999999
$ws.Range("I3").Value = 0.8824726436021215
$ws.Range("J3").Value = 0.8824726436021214
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.831039666666666
$ws.Range("N3").Value = 23.493119
$ws.Range("O3").Value = 0.3944590147859676
$ws.Range("P3").Value = 0.3944590147859677
$ws.Range("Q3").Value = 184.8300411173348
$ws.Range("R3").Value = 1663.470370056014
$ws.Range("S3").Value = 0.3480992895708612
$ws.Range("T3").Value = 0.3480992895708612
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 23.60223533333333
$ws.Range("H4").Value = 70.80670599999999
$ws.Range("I4").Value = 0.8824726436021215
$ws.Range("J4").Value = 0.8824726436021214
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.804271333333333
$ws.Range("N4").Value = 5.412813999999999
$ws.Range("O4").Value = 0.09088334663693197
$ws.Range("P4").Value = 0.09088334663693197
$ws.Range("Q4").Value = 42.58483661452043
$ws.Range("R4").Value = 383.2635295306839
$ws.Range("S4").Value = 0.08020206716610133
$ws.Range("T4").Value = 0.08020206716610133
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 23.60223533333333
$ws.Range("H5").Value = 70.80670599999999
$ws.Range("I5").Value = 0.8824726436021215
$ws.Range("J5").Value = 0.8824726436021214
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.220799666666667
$ws.Range("N5").Value = 9.662399000000001
$ws.Range("O5").Value = 0.1622356056685755
$ws.Range("P5").Value = 0.1622356056685755
$ws.Range("Q5").Value = 76.01807169418821
$ws.Range("R5").Value = 684.1626452476939
$ws.Range("S5").Value = 0.1431684838207392
$ws.Range("T5").Value = 0.1431684838207392
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 23.60223533333333
$ws.Range("H6").Value = 70.80670599999999
$ws.Range("I6").Value = 0.8824726436021215
$ws.Range("J6").Value = 0.8824726436021214
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.498796666666667
$ws.Range("N6").Value = 7.49639
$ws.Range("O6").Value = 0.1258674343688201
$ws.Range("P6").Value = 0.1258674343688201
$ws.Range("Q6").Value = 58.97718697681555
$ws.Range("R6").Value = 530.7946827913399
$ws.Range("S6").Value = 0.1110745675508692
$ws.Range("T6").Value = 0.1110745675508692
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.633202
$ws.Range("H7").Value = 7.899606
$ws.Range("I7").Value = 0.09845375648791208
$ws.Range("J7").Value = 0.09845375648791205
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.497699333333333
$ws.Range("N7").Value = 13.493098
$ws.Range("O7").Value = 0.2265545985397048
$ws.Range("P7").Value = 0.2265545985397048
$ws.Range("Q7").Value = 11.843350879932
$ws.Range("R7").Value = 106.590157919388
$ws.Range("S7").Value = 0.02230515127584478
$ws.Range("T7").Value = 0.02230515127584477
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.633202
$ws.Range("H8").Value = 7.899606
$ws.Range("I8").Value = 0.09845375648791208
$ws.Range("J8").Value = 0.09845375648791205
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.831039666666666
$ws.Range("N8").Value = 23.493119
$ws.Range("O8").Value = 0.3944590147859676
$ws.Range("P8").Value = 0.3944590147859677
$ws.Range("Q8").Value = 20.620709312346
$ws.Range("R8").Value = 185.586383811114
$ws.Range("S8").Value = 0.03883597178619937
$ws.Range("T8").Value = 0.03883597178619937
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.633202
$ws.Range("H9").Value = 7.899606
$ws.Range("I9").Value = 0.09845375648791208
$ws.Range("J9").Value = 0.09845375648791205
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.804271333333333
$ws.Range("N9").Value = 5.412813999999999
$ws.Range("O9").Value = 0.09088334663693197
$ws.Range("P9").Value = 0.09088334663693197
$ws.Range("Q9").Value = 4.751010883476
$ws.Range("R9").Value = 42.759097951284
$ws.Range("S9").Value = 0.008947806878599003
$ws.Range("T9").Value = 0.008947806878599
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.633202
$ws.Range("H10").Value = 7.899606
$ws.Range("I10").Value = 0.09845375648791208
$ws.Range("J10").Value = 0.09845375648791205
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.220799666666667
$ws.Range("N10").Value = 9.662399000000001
$ws.Range("O10").Value = 0.1622356056685755
$ws.Range("P10").Value = 0.1622356056685755
$ws.Range("Q10").Value = 8.481016123866
$ws.Range("R10").Value = 76.329145114794
$ws.Range("S10").Value = 0.01597270481416286
$ws.Range("T10").Value = 0.01597270481416286
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.633202
$ws.Range("H11").Value = 7.899606
$ws.Range("I11").Value = 0.09845375648791208
$ws.Range("J11").Value = 0.09845375648791205
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.498796666666667
$ws.Range("N11").Value = 7.49639
$ws.Range("O11").Value = 0.1258674343688201
$ws.Range("P11").Value = 0.1258674343688201
$ws.Range("Q11").Value = 6.579836380260001
$ws.Range("R11").Value = 59.21852742234
$ws.Range("S11").Value = 0.01239212173310607
$ws.Range("T11").Value = 0.01239212173310606
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.5101343333333334
$ws.Range("H12").Value = 1.530403
$ws.Range("I12").Value = 0.01907359990996641
$ws.Range("J12").Value = 0.0190735999099664
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.497699333333333
$ws.Range("N12").Value = 13.493098
$ws.Range("O12").Value = 0.2265545985397048
$ws.Range("P12").Value = 0.2265545985397048
$ws.Range("Q12").Value = 2.294430850943778
$ws.Range("R12").Value = 20.649877658494
$ws.Range("S12").Value = 0.004321211770309388
$ws.Range("T12").Value = 0.004321211770309388
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.5101343333333334
$ws.Range("H13").Value = 1.530403
$ws.Range("I13").Value = 0.01907359990996641
$ws.Range("J13").Value = 0.0190735999099664
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 7.831039666666666
$ws.Range("N13").Value = 23.493119
$ws.Range("O13").Value = 0.3944590147859676
$ws.Range("P13").Value = 0.3944590147859677
$ws.Range("Q13").Value = 3.994882199661889
$ws.Range("R13").Value = 35.953939796957
$ws.Range("S13").Value = 0.00752375342890707
$ws.Range("T13").Value = 0.007523753428907069
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.5101343333333334
$ws.Range("H14").Value = 1.530403
$ws.Range("I14").Value = 0.01907359990996641
$ws.Range("J14").Value = 0.0190735999099664
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.804271333333333
$ws.Range("N14").Value = 5.412813999999999
$ws.Range("O14").Value = 0.09088334663693197
$ws.Range("P14").Value = 0.09088334663693197
$ws.Range("Q14").Value = 0.9204207537824444
$ws.Range("R14").Value = 8.283786784041999
$ws.Range("S14").Value = 0.001733472592231631
$ws.Range("T14").Value = 0.001733472592231631
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.5101343333333334
$ws.Range("H15").Value = 1.530403
$ws.Range("I15").Value = 0.01907359990996641
$ws.Range("J15").Value = 0.0190735999099664
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.220799666666667
$ws.Range("N15").Value = 9.662399000000001
$ws.Range("O15").Value = 0.1622356056685755
$ws.Range("P15").Value = 0.1622356056685755
$ws.Range("Q15").Value = 1.643040490755222
$ws.Range("R15").Value = 14.787364416797
$ws.Range("S15").Value = 0.003094417033673488
$ws.Range("T15").Value = 0.003094417033673487
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.5101343333333334
$ws.Range("H16").Value = 1.530403
$ws.Range("I16").Value = 0.01907359990996641
$ws.Range("J16").Value = 0.0190735999099664
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 2.498796666666667
$ws.Range("N16").Value = 7.49639
$ws.Range("O16").Value = 0.1258674343688201
$ws.Range("P16").Value = 0.1258674343688201
$ws.Range("Q16").Value = 1.274721971685556
$ws.Range("R16").Value = 11.47249774517
$ws.Range("S16").Value = 0.00240074508484483
$ws.Range("T16").Value = 0.002400745084844829
